$wb = $excel.ActiveWorkbook

# --- workbook.xml: active tab moves from network_weights (5) to optimization_parameters (6) ---
# (handled implicitly below by activating network_weights first, then optimization_parameters last,
#  so optimization_parameters ends up the ActiveSheet / activeTab, and tabSelected="1" ends on it.)

$wsNetworkWeights = $wb.Worksheets.Item("network_weights")
$wsOptParams = $wb.Worksheets.Item("optimization_parameters")

# --- sheet6 (network_weights): selection moves from C7 to E49, tabSelected removed ---
$wsNetworkWeights.Activate()
$wsNetworkWeights.Range("E49").Select()

# --- sheet7 (optimization_parameters): row 1 trimmed to A:B only ---
$wsOptParams.Range("C1:F1").ClearContents()

# Row 8: "Model"/"Sigmoid" -> "production_function"/"Sigmoid" (style to match header row)
$wsOptParams.Range("A8").Value = "production_function"
$wsOptParams.Range("A8").Font.Name = "Verdana"

# Remove old "Deletion" row (was row 16: Deletion/0/3)
$wsOptParams.Rows.Item(16).Delete()

# Insert new row 9 for "L_curve" / 1
$wsOptParams.Rows.Item(9).Insert()
$wsOptParams.Range("A9").Value = "L_curve"
$wsOptParams.Range("A9").Font.Name = "Verdana"
$wsOptParams.Range("B9").Value = 1
$wsOptParams.Range("B9").NumberFormat = "0.00E+00"

# --- tab / selection: optimization_parameters becomes active, whole row 17 selected ---
$wsOptParams.Activate()
$wsOptParams.Rows.Item(17).Select()

Write-Output "done"
